$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medium traffic density")
$ws3 = $wb.Worksheets.Item("High traffic density")

# Row 25 scenario label changes from the "D2C2 alt" outlier string to the regular "D2C2" label
$ws.Range("B25").Value = "D2C2"

# Updated per-run measurements for the medium-traffic-density dataset (rows 17-26)
$ws.Range("A17").Value = 1
$ws.Range("C17").Value = 158
$ws.Range("D17").Value = 101
$ws.Range("E17").Value = 4169.6757416205901
$ws.Range("F17").Value = 546.12821782178196
$ws.Range("H17").Value = 0.91120313498790995
$ws.Range("I17").Value = 59
$ws.Range("J17").Value = 7.1016949152542397
$ws.Range("K17").Value = 67
$ws.Range("L17").Value = 19

$ws.Range("A18").Value = 2
$ws.Range("C18").Value = 148
$ws.Range("D18").Value = 85
$ws.Range("E18").Value = 3666.1062516249399
$ws.Range("F18").Value = 500.32352941176498
$ws.Range("H18").Value = 0.75650790993554895
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 7.5833332999999996
$ws.Range("K18").Value = 45
$ws.Range("L18").Value = 14

$ws.Range("A19").Value = 3
$ws.Range("C19").Value = 182
$ws.Range("D19").Value = 98
$ws.Range("E19").Value = 3715.1805563173498
$ws.Range("F19").Value = 496.375
$ws.Range("H19").Value = 0.87801800158121901
$ws.Range("I19").Value = 66
$ws.Range("J19").Value = 8
$ws.Range("K19").Value = 70
$ws.Range("L19").Value = 28

$ws.Range("A20").Value = 4
$ws.Range("C20").Value = 148
$ws.Range("D20").Value = 96
$ws.Range("E20").Value = 3753.1046722368801
$ws.Range("F20").Value = 504.91354166669998
$ws.Range("H20").Value = 0.79966317231164896
$ws.Range("I20").Value = 42
$ws.Range("J20").Value = 6.3333332999999996
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 18

$ws.Range("A21").Value = 5
$ws.Range("C21").Value = 181
$ws.Range("D21").Value = 98
$ws.Range("E21").Value = 3695.9154167349002
$ws.Range("F21").Value = 481.52040816326502
$ws.Range("H21").Value = 0.84160326539679098
$ws.Range("I21").Value = 49
$ws.Range("J21").Value = 7.83673469387755
$ws.Range("K21").Value = 62
$ws.Range("L21").Value = 19

$ws.Range("A22").Value = 6
$ws.Range("C22").Value = 171
$ws.Range("D22").Value = 87
$ws.Range("E22").Value = 3861.4096512134502
$ws.Range("F22").Value = 534.59367816092004
$ws.Range("H22").Value = 0.81070380785793195
$ws.Range("I22").Value = 62
$ws.Range("J22").Value = 6.1935483870967696
$ws.Range("K22").Value = 70
$ws.Range("L22").Value = 24

$ws.Range("A23").Value = 7
$ws.Range("C23").Value = 165
$ws.Range("D23").Value = 101
$ws.Range("E23").Value = 3765.5085450500001
$ws.Range("F23").Value = 532.55841584158395
$ws.Range("H23").Value = 0.87573879999048099
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 5.8085106382978697
$ws.Range("K23").Value = 51
$ws.Range("L23").Value = 21

$ws.Range("A24").Value = 8
$ws.Range("C24").Value = 155
$ws.Range("D24").Value = 95
$ws.Range("E24").Value = 4002.6993699562099
$ws.Range("F24").Value = 550.66105263157897
$ws.Range("H24").Value = 0.88673982697217202
$ws.Range("I24").Value = 57
$ws.Range("J24").Value = 5.59649122807018
$ws.Range("K24").Value = 74
$ws.Range("L24").Value = 24

$ws.Range("A25").Value = 9
$ws.Range("C25").Value = 171
$ws.Range("D25").Value = 105
$ws.Range("E25").Value = 3797.5452887718102
$ws.Range("F25").Value = 532.24190476190495
$ws.Range("H25").Value = 0.98544205208833802
$ws.Range("I25").Value = 58
$ws.Range("J25").Value = 5.5344827586206904
$ws.Range("K25").Value = 66
$ws.Range("L25").Value = 19

$ws.Range("A26").Value = 10
$ws.Range("C26").Value = 172
$ws.Range("D26").Value = 101
$ws.Range("E26").Value = 3626.4529817215798
$ws.Range("F26").Value = 487.112871287129
$ws.Range("H26").Value = 0.85595036482068398
$ws.Range("I26").Value = 51
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = 48
$ws.Range("L26").Value = 21

# Drop the old 11th run (previously row 27); Average/Std-dev rows shift up to rows 27/28
$ws.Rows.Item(27).Delete()

# Re-establish the shared ratio formula across the (now 10-row) sample G17:G26
$ws.Range("G17:G26").Formula = "=E17/F17"

# Rebuild the Average / Standard-dev summary rows against the new G17:G26 sample range
foreach ($col in @("C","D","E","F","G","H","I","J","K","L")) {
    $ws.Range("$col`27").Formula = "=AVERAGE($col`17:$col`26)"
}
foreach ($col in @("C","D","E","F","G","H","I","J","K","L")) {
    $ws.Range("$col`28").Formula = "=_xlfn.STDEV.P($col`17:$col`26)"
}

# Medium traffic density becomes the active/visible sheet (was High traffic density)
$ws.Activate()
